$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset for rows 2-29 (A: timestamp, B: label, C-H: sensor axes)
$data = @(
    @(0, "falling", -0.7498788833618164, 0.7729501724243164, -1.214103698730469, -0.0186313893646001, -0.1148426681756973, 0.057115901261568),
    @(100, "falling", -0.8813939094543457, 0.8021388053894043, -1.258926272392273, 0.0178678091615438, -0.1070541366934776, 0.1860084682703018),
    @(200, "falling", -0.6222906112670898, 0.7483844757080078, -1.280380129814148, -0.007941247895359899, 0.0074830991216003, 0.0245873257517814),
    @(300, "falling", -0.7319130897521973, 0.7634215354919434, -1.176308631896973, -0.039248090237379, -0.008552113547921099, 0.0209221355617046),
    @(400, "falling", -0.568336009979248, 0.7112784385681152, -1.411303043365478, 0.00167987938039, -0.0282525178045034, 0.0545197241008281),
    @(500, "falling", -0.4251332283020019, 0.7228684425354004, -1.317938923835754, 0.00167987938039, -0.0485637858510017, 0.0387899428606033),
    @(600, "falling", -0.824821949005127, 0.7479877471923828, -1.330062508583069, 0.0103847095742821, 0.0097738439217209, -0.0050396383740007),
    @(700, "falling", -0.9185628890991212, 0.8159165382385254, -1.254562854766846, -0.0229074470698833, 0.0265726372599601, 0.0036651915870606),
    @(800, "falling", -0.8128528594970703, 0.7384524345397949, -1.323878526687622, 0.011148290708661, 0.0045814891345798, -0.0197004042565822),
    @(900, "falling", -0.5396947860717773, 0.4983515739440918, -2.386690616607666, 0.0502436682581901, -0.0682641938328743, 0.0099265603348612),
    @(1000, "falling", -0.1267647743225097, 0.3085379600524902, -2.642498016357422, -0.0494800843298435, -0.3627012372016907, 0.144622340798378),
    @(1100, "falling", -0.2359266281127929, -0.2183656692504882, -3.997431516647339, -0.0519235469400882, -0.6624833345413208, 0.2745839357376098),
    @(1200, "falling", -0.7616920471191406, -1.078659534454346, -5.043187141418457, 0.4827362596988678, -0.633009135723114, 0.4114177525043487),
    @(1300, "falling", -2.656810522079468, -2.278046607971191, -4.701850414276123, 0.427911102771759, 0.3917173445224762, 0.2205223590135574),
    @(1400, "falling", -4.357179164886475, -1.795724391937256, -4.692870140075684, 0.1753183305263519, 0.8854491710662842, -0.3863722681999206),
    @(1500, "falling", -2.303323268890381, -1.950790405273437, -3.430110692977905, 0.0439822971820831, 1.484707951545715, -1.550223231315613),
    @(1600, "falling", -4.76627254486084, -2.148142337799072, -2.799470901489258, -0.4473060667514801, -2.502867698669434, -7.500508785247803),
    @(1700, "falling", -5.706301212310791, -2.392673492431641, 3.553245782852173, -1.311222195625305, -0.5021312236785889, -1.749059915542602),
    @(1800, "falling", 3.05734920501709, -6.071773529052734, 0.5925030708312988, -0.0354301854968071, 0.5714644193649292, 0.3466660380363464),
    @(1900, "falling", -0.7244548797607422, 0.7260329723358154, -0.8456218242645264, 0.0313068442046642, 1.005331516265869, 1.156062483787537),
    @(2000, "falling", -0.8577280044555664, 1.746235489845276, -1.253829836845398, 0.0287106670439243, 0.3336851298809051, 0.0453567430377006),
    @(2100, "falling", -0.9032430648803712, 0.7269128561019897, -1.814414024353028, 0.0639881342649459, 0.8878926634788513, 0.5963571667671204),
    @(2200, "falling", 0.3531332015991211, 3.250727415084839, -0.5509147644042969, 0.1010981947183609, -0.2167044430971145, -0.2678644061088562),
    @(2300, "falling", -0.1614398956298828, 0.2837758064270019, -2.02338171005249, 0.08964447677135461, -0.2109012305736541, 0.0058032199740409),
    @(2400, "falling", -0.5078344345092773, 0.6439783573150635, -1.895825386047364, -0.0484110713005065, 0.1492038369178772, 0.1434006094932556),
    @(2500, "falling", -0.4796314239501953, 0.7996485829353333, -0.4695221781730652, -0.1841758787631988, 0.1082758679986, 0.0047342055477201),
    @(2600, "falling", -0.5036382675170898, 1.3280930519104, -1.352860689163208, -0.102472648024559, -0.0048869219608604, 0.0861319974064827),
    @(2700, "falling", -0.525787353515625, 1.078789949417114, -1.021092414855957, -0.0229074470698833, 0.0600175112485885, 0.0545197241008281)
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $row[$c - 1]
    }
    $r++
}
